$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'257.11"
$ws.Range("E2").Value = "'-0.18%"
$ws.Range("D3").Value = "'26.97"
$ws.Range("E3").Value = "'-0.98%"
$ws.Range("E4").Value = "'-11.52%"
$ws.Range("D6").Value = "'6.624"
$ws.Range("E6").Value = "'-0.69%"
$ws.Range("E7").Value = "'-0.60%"
$ws.Range("D8").Value = "'0.9423"
$ws.Range("E8").Value = "'-7.86%"
$ws.Range("E9").Value = "'-0.59%"
$ws.Range("D10").Value = "'0.04158"
$ws.Range("E10").Value = "'16.49%"
$ws.Range("D11").Value = "'0.07099"
$ws.Range("E11").Value = "'-1.31%"
$ws.Range("D12").Value = "'0.03185"
$ws.Range("E12").Value = "'1.59%"
$ws.Range("D13").Value = "'0.09160"
$ws.Range("E13").Value = "'-0.83%"
$ws.Range("D14").Value = "'0.001539"
$ws.Range("E14").Value = "'-0.38%"
$ws.Range("D15").Value = "'0.0006062"
$ws.Range("E15").Value = "'-0.25%"
$ws.Range("D16").Value = "'0.006231"
$ws.Range("E16").Value = "'9.67%"
$ws.Range("D17").Value = "'3.523"
$ws.Range("E17").Value = "'0.74%"
$ws.Range("E19").Value = "'-0.02%"
$ws.Range("D20").Value = "'0.3053"
$ws.Range("E20").Value = "'-2.88%"
$ws.Range("E21").Value = "'-0.45%"
$ws.Range("D22").Value = "'3.832"
$ws.Range("E22").Value = "'8.84%"
$ws.Range("D23").Value = "'0.04227"
$ws.Range("E23").Value = "'1.19%"
$ws.Range("D24").Value = "'0.001226"
$ws.Range("E24").Value = "'0.67%"
$ws.Range("D25").Value = "'0.004290"
$ws.Range("E25").Value = "'-4.99%"
$ws.Range("E26").Value = "'0.07%"
$ws.Range("D27").Value = "'0.0001938"
$ws.Range("E27").Value = "'30.57%"
$ws.Range("D40").Value = "'0.03827"
$ws.Range("E40").Value = "'0.31%"
$ws.Range("D41").Value = "'0.006210"
$ws.Range("E41").Value = "'-5.55%"
$ws.Range("E42").Value = "'-0.06%"
$ws.Range("D43").Value = "'0.002430"
$ws.Range("E43").Value = "'10.53%"
$ws.Range("D44").Value = "'0.01144"
$ws.Range("E44").Value = "'6.24%"
$ws.Range("D45").Value = "'0.00005464"
$ws.Range("E45").Value = "'0.67%"
$ws.Range("E46").Value = "'0.10%"
$ws.Range("D47").Value = "'0.05123"
$ws.Range("D48").Value = "'0.2347"
$ws.Range("E48").Value = "'10,402.48%"
$ws.Range("D49").Value = "'0.00002101"
$ws.Range("E49").Value = "'0.10%"
$ws.Range("D50").Value = "'0.0002001"
$ws.Range("E50").Value = "'0.10%"
